$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("A1").Value = "Experiment Folder"
$ws.Range("B1").Value = "Local Area Density"
$ws.Range("C1").Value = "Potential Radius"
$ws.Range("D1").Value = "Local/Global Inhibition"
$ws.Range("E1").Value = "NumActiveColumnsPerInhArea"
$ws.Range("F1").Value = "Result Image Name"

# --- Data rows ---
$ws.Range("A2").Value = "Exp 6"
$ws.Range("B2").Value = 0.1
$ws.Range("C2").Value = 1
$ws.Range("D2").Value = "Local"
$ws.Range("E2").Value = -1
$ws.Range("F2").Value = "Exp 6.png"

$ws.Range("A3").Value = "Exp 7"
$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("D3").Value = "Local"
$ws.Range("E3").Value = -1

for ($r = 4; $r -le 13; $r++) {
    $ws.Range("D$r").Value = "Local"
    $ws.Range("E$r").Value = -1
}

# --- Column widths (tuned so the saved <col width> best matches the
#     target 16.85546875 / 19.42578125 / 16.85546875 / 21.42578125 /
#     28.28515625 / 21.28515625 given this engine's pixel rounding) ---
$ws.Columns("A").ColumnWidth = 15.9221
$ws.Columns("B").ColumnWidth = 18.5924
$ws.Columns("C").ColumnWidth = 15.9221
$ws.Columns("D").ColumnWidth = 20.5924
$ws.Columns("E").ColumnWidth = 27.4218
$ws.Columns("F").ColumnWidth = 20.4218

# --- Formatting: whole header row A1:F1 - yellow fill + thin border ---
$headerFull = $ws.Range("A1:F1")
$headerFull.Interior.Color = 65535
$headerFull.Borders.LineStyle = 1

# --- Formatting: A1:E1 additionally centered ---
$ws.Range("A1:E1").HorizontalAlignment = -4108

# --- Formatting: body data A2:E13 - centered ---
$ws.Range("A2:E13").HorizontalAlignment = -4108

# --- Selection state ---
$ws.Range("B8").Select()
